$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 from "N" to "Y" (Data Driven Extended With Docker)
$ws.Range("B3").Value = "Y"

# Move the active selection to B3
$ws.Range("B3").Select()
